$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4488.25
$ws.Range("J64").Value = 3654.3333
$ws.Range("L64").Value = 3654.3333
$ws.Range("N64").Value = -4150.3333

$ws.Range("H67").Value = 4488.25
$ws.Range("J67").Value = 3654.3333
$ws.Range("L67").Value = 3654.3333
$ws.Range("N67").Value = -5370.3333

$ws.Range("H92").Value = 1794.5333
$ws.Range("I92").Value = 1751.5
$ws.Range("K92").Value = 1751.5
$ws.Range("M92").Value = -503.5

$ws.Range("H107").Value = 4020.6667
$ws.Range("I107").Value = 3814.8
$ws.Range("J107").Value = 5050
$ws.Range("K107").Value = 3814.8
$ws.Range("L107").Value = 5050
$ws.Range("M107").Value = -1894.8
$ws.Range("N107").Value = -8890

$ws.Range("H111").Value = 3288.1667
$ws.Range("I111").Value = 3288.1667
$ws.Range("K111").Value = 9864.500100000001
$ws.Range("M111").Value = -6797.500100000001

$ws.Range("H132").Value = 7759376.5
$ws.Range("I132").Value = 11116447
$ws.Range("K132").Value = 33349341
$ws.Range("M132").Value = -33346811

$ws.Range("H137").Value = 1643.8422
$ws.Range("I137").Value = 988
$ws.Range("K137").Value = 2964
$ws.Range("M137").Value = -414

$ws.Range("H138").Value = 519327.06
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 519327.06
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 1557981.18
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -1568261.18

$ws.Range("H141").Value = 656.8182
$ws.Range("I141").Value = 656.8182
$ws.Range("K141").Value = 1970.4546
$ws.Range("M141").Value = 3209.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5768.4653
$ws.Range("I32").Value = 5379.8423
$ws.Range("K32").Value = 5379.8423
$ws.Range("M32").Value = -5092.8423

$ws.Range("H74").Value = 1674
$ws.Range("I74").Value = 887.53845
$ws.Range("K74").Value = 887.53845
$ws.Range("M74").Value = -13.53845000000001

$ws.Range("H77").Value = 1674
$ws.Range("I77").Value = 887.53845
$ws.Range("K77").Value = 4437.69225
$ws.Range("M77").Value = -69.69225000000006

$ws.Range("H110").Value = 1336.2106
$ws.Range("I110").Value = 1065.6666
$ws.Range("K110").Value = 1065.6666
$ws.Range("M110").Value = 979.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 40445
$ws.Range("I51").Value = 30000
$ws.Range("J51").Value = 43926.668
$ws.Range("K51").Value = 30000
$ws.Range("L51").Value = 43926.668
$ws.Range("M51").Value = -29509
$ws.Range("N51").Value = -44908.668

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H94").Value = 16667516
$ws.Range("I94").Value = 17857894
$ws.Range("K94").Value = 17857894
$ws.Range("M94").Value = -17857443

$ws.Range("H129").Value = 49759.332
$ws.Range("J129").Value = 49759.332
$ws.Range("L129").Value = 49759.332
$ws.Range("N129").Value = -59759.332

$ws.Range("H134").Value = 4553.324
$ws.Range("I134").Value = 1122.2413
$ws.Range("K134").Value = 3366.7239
$ws.Range("M134").Value = -831.7239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2326.6667
$ws.Range("I31").Value = 2490
$ws.Range("K31").Value = 2490
$ws.Range("M31").Value = -2195

$ws.Range("H34").Value = 2326.6667
$ws.Range("I34").Value = 2490
$ws.Range("K34").Value = 2490
$ws.Range("M34").Value = -2288

$ws.Range("H132").Value = 1738.9512
$ws.Range("I132").Value = 1327.4242
$ws.Range("J132").Value = 3436.5
$ws.Range("K132").Value = 3982.2726
$ws.Range("L132").Value = 10309.5
$ws.Range("M132").Value = -1452.2726
$ws.Range("N132").Value = -15369.5

$ws.Range("H134").Value = 830.8214
$ws.Range("I134").Value = 711.4783
$ws.Range("K134").Value = 2134.4349
$ws.Range("M134").Value = 400.5650999999998

$ws.Range("H138").Value = 111417.664
$ws.Range("J138").Value = 111417.664
$ws.Range("L138").Value = 111417.664
$ws.Range("N138").Value = -121697.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6359.7085
$ws.Range("I56").Value = 6359.7085
$ws.Range("K56").Value = 6359.7085
$ws.Range("M56").Value = -5829.7085

$ws.Range("H119").Value = 6824.1816

$ws.Range("H120").Value = 8766
$ws.Range("J120").Value = 11999
$ws.Range("L120").Value = 35997
$ws.Range("N120").Value = -45673

$ws.Range("H121").Value = 376.66666
$ws.Range("I121").Value = 376.66666
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 1129.99998
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 180.0000199999999
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 811.3889
$ws.Range("I122").Value = 492.22223
$ws.Range("J122").Value = 1130.5555
$ws.Range("K122").Value = 4430.00007
$ws.Range("L122").Value = 10174.9995
$ws.Range("M122").Value = -1980.00007
$ws.Range("N122").Value = -15074.9995

$ws.Range("H123").Value = 2998
$ws.Range("I123").Value = 2990
$ws.Range("J123").Value = 2998.7273
$ws.Range("K123").Value = 8970
$ws.Range("L123").Value = 8996.1819
$ws.Range("M123").Value = -6520
$ws.Range("N123").Value = -13896.1819

$ws.Range("H126").Value = 5625.294
$ws.Range("I126").Value = 2515
$ws.Range("J126").Value = 6040
$ws.Range("K126").Value = 7545
$ws.Range("L126").Value = 18120
$ws.Range("M126").Value = -2605
$ws.Range("N126").Value = -28000

$ws.Range("H131").Value = 23257004
$ws.Range("J131").Value = 1518.9375
$ws.Range("L131").Value = 4556.8125
$ws.Range("N131").Value = -14636.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 798.53845
$ws.Range("I97").Value = 796.6667
$ws.Range("J97").Value = 802.75
$ws.Range("K97").Value = 796.6667
$ws.Range("L97").Value = 802.75
$ws.Range("M97").Value = -300.6667
$ws.Range("N97").Value = -1794.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1525
$ws.Range("J46").Value = 1700
$ws.Range("L46").Value = 1700
$ws.Range("N46").Value = -2076

$ws.Range("H61").Value = 2199.111
$ws.Range("I61").Value = 1746.75
$ws.Range("J61").Value = 2561
$ws.Range("K61").Value = 1746.75
$ws.Range("L61").Value = 2561
$ws.Range("M61").Value = -1544.75
$ws.Range("N61").Value = -2965

$ws.Range("H93").Value = 919.6
$ws.Range("I93").Value = 919.6
$ws.Range("K93").Value = 919.6
$ws.Range("M93").Value = 328.4

$ws.Range("H109").Value = 30285
$ws.Range("J109").Value = 30285
$ws.Range("L109").Value = 30285
$ws.Range("N109").Value = -33059

$ws.Range("H113").Value = 2199.111
$ws.Range("I113").Value = 1746.75
$ws.Range("J113").Value = 2561
$ws.Range("K113").Value = 1746.75
$ws.Range("L113").Value = 2561
$ws.Range("M113").Value = 423.25
$ws.Range("N113").Value = -6901

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H136").Value = 1463.8096
$ws.Range("I136").Value = 1276.0667
$ws.Range("J136").Value = 1933.1666
$ws.Range("K136").Value = 3828.2001
$ws.Range("L136").Value = 5799.4998
$ws.Range("M136").Value = -1278.2001
$ws.Range("N136").Value = -10899.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 853.1429000000001
$ws.Range("I100").Value = 989.44446
$ws.Range("K100").Value = 1978.88892
$ws.Range("M100").Value = -1437.88892

$ws.Range("H132").Value = 2624.7856
$ws.Range("I132").Value = 2524.44
$ws.Range("J132").Value = 3461
$ws.Range("K132").Value = 7573.32
$ws.Range("L132").Value = 10383
$ws.Range("M132").Value = -5043.32
$ws.Range("N132").Value = -15443

$ws.Range("H136").Value = 514
$ws.Range("I136").Value = 266.25
$ws.Range("J136").Value = 2000.5
$ws.Range("K136").Value = 798.75
$ws.Range("L136").Value = 6001.5
$ws.Range("M136").Value = 1751.25
$ws.Range("N136").Value = -11101.5
